# Menambahkan sebuah komentar pada masing masing fitur agar dapat mudah di pahami
$wb = $excel.ActiveWorkbook

# Update remaining quota ("Kuota") for scholarship B01001 on sheet "Beasiswa"
# one seat was just awarded, so the available quota drops from 200 to 199
$wsBeasiswa = $wb.Worksheets.Item("Beasiswa")
$wsBeasiswa.Range("F2").Value = 199

# Record the new scholarship award on sheet "Pemberian"
$wsPemberian = $wb.Worksheets.Item("Pemberian")

# NISN (A) must stay text so leading zeros are preserved
$wsPemberian.Cells.Item(3, 1).NumberFormat = "@"
$wsPemberian.Cells.Item(3, 1).Value = "0012345678"
$wsPemberian.Cells.Item(3, 1).Style = "Normal"

$wsPemberian.Cells.Item(3, 2).Value = "B01001"

# Tanggal (C) must stay text in yyyy-mm-dd form rather than becoming a date serial
$wsPemberian.Cells.Item(3, 3).NumberFormat = "@"
$wsPemberian.Cells.Item(3, 3).Value = "2025-12-01"
$wsPemberian.Cells.Item(3, 3).Style = "Normal"
